# Apply updates described in the diff:
# - Fill in previously-missing values for rows 199-200
# - Append new data rows 201-204
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 199 ----
$ws.Range("B199").Value = 321.5

# ---- Row 200 ----
$ws.Range("B200").Value = 322.132
$ws.Range("C200").Value = 5
$ws.Range("D200").Value = 2.38054202
$ws.Range("E200").Value = 0.02637591034548059
$ws.Range("I200").Value = 0.7555555555555555
$ws.Range("W200").Value = 0.02623864148733057
$ws.Range("X200").Value = 0.02455862441701219
$ws.Range("Y200").Value = 0.02599613875449484

# ---- Row 201 (new) ----
$ws.Range("A201").Value = 45900
$ws.Range("B201").Value = 323.364
$ws.Range("C201").Value = 4.5
$ws.Range("D201").Value = 2.78904728
$ws.Range("E201").Value = 0.0269515351998173
$ws.Range("F201").Value = 5
$ws.Range("G201").Value = 0.7555555555555555
$ws.Range("H201").Value = 0.02637591034548059
$ws.Range("I201").Value = 0.8146953405017922
$ws.Range("J201").Value = 0.9185185185185186
$ws.Range("K201").Value = 1.170250896057348
$ws.Range("L201").Value = 1.117037037037037
$ws.Range("M201").Value = 0.8412186379928316
$ws.Range("N201").Value = 0.6595238095238095
$ws.Range("O201").Value = 0.4892473118279569
$ws.Range("P201").Value = 5
$ws.Range("Q201").Value = 6.6
$ws.Range("R201").Value = 6.5
$ws.Range("S201").Value = 5
$ws.Range("T201").Value = 4.3
$ws.Range("U201").Value = 3.3
$ws.Range("V201").Value = 0.02638692659292648
$ws.Range("W201").Value = 0.02819815360789871
$ws.Range("X201").Value = 0.0280586362304498
$ws.Range("Y201").Value = 0.02720054445549156

# ---- Row 202 (new) ----
$ws.Range("A202").Value = 45930
$ws.Range("B202").Value = 324.368
$ws.Range("C202").Value = 4.8
$ws.Range("D202").Value = 2.68936045
$ws.Range("E202").Value = 0.02896852733478372
$ws.Range("F202").Value = 4.5
$ws.Range("G202").Value = 0.8146953405017922
$ws.Range("H202").Value = 0.0269515351998173
$ws.Range("I202").Value = 0.9248148148148148
$ws.Range("J202").Value = 0.7555555555555555
$ws.Range("K202").Value = 0.9185185185185186
$ws.Range("L202").Value = 1.170250896057348
$ws.Range("M202").Value = 1.117037037037037
$ws.Range("N202").Value = 0.8412186379928316
$ws.Range("O202").Value = 0.6595238095238095
$ws.Range("P202").Value = 4.5
$ws.Range("Q202").Value = 5
$ws.Range("R202").Value = 6.6
$ws.Range("S202").Value = 6.5
$ws.Range("T202").Value = 5
$ws.Range("U202").Value = 4.3
$ws.Range("V202").Value = 0.02695254058730978
$ws.Range("W202").Value = 0.02914323536072327
$ws.Range("X202").Value = 0.0283830056392974
$ws.Range("Y202").Value = 0.02893361365055446

# ---- Row 203 (new) ----
$ws.Range("A203").Value = 45961
# B203 intentionally left blank (missing data)
$ws.Range("C203").Value = 4.7
$ws.Range("D203").Value = 2.79604823
$ws.Range("E203").Value = 0.02977916268051928
$ws.Range("F203").Value = 4.8
$ws.Range("G203").Value = 0.9248148148148148
$ws.Range("H203").Value = 0.02896852733478372
$ws.Range("I203").Value = 0.8086021505376344
$ws.Range("J203").Value = 0.8146953405017922
$ws.Range("K203").Value = 0.7555555555555555
$ws.Range("L203").Value = 0.9185185185185186
$ws.Range("M203").Value = 1.170250896057348
$ws.Range("N203").Value = 1.117037037037037
$ws.Range("O203").Value = 0.8412186379928316
$ws.Range("P203").Value = 4.8
$ws.Range("Q203").Value = 4.5
$ws.Range("R203").Value = 5
$ws.Range("S203").Value = 6.6
$ws.Range("T203").Value = 6.5
$ws.Range("U203").Value = 5
$ws.Range("V203").Value = 0.02893488133016218
$ws.Range("W203").Value = 0.03087154223420583
$ws.Range("X203").Value = 0.0304749717201673
$ws.Range("Y203").Value = 0.03002951986028322

# ---- Row 204 (new) ----
$ws.Range("A204").Value = 45991
# B204 intentionally left blank (missing data)
$ws.Range("C204").Value = 4.6
$ws.Range("D204").Value = 2.74488063
# E204 intentionally left blank (missing data)
$ws.Range("F204").Value = 4.7
$ws.Range("G204").Value = 0.8086021505376344
$ws.Range("H204").Value = 0.02977916268051928
# I204 intentionally left blank (missing data)
$ws.Range("J204").Value = 0.9248148148148148
$ws.Range("K204").Value = 0.8146953405017922
$ws.Range("L204").Value = 0.7555555555555555
$ws.Range("M204").Value = 0.9185185185185186
$ws.Range("N204").Value = 1.170250896057348
$ws.Range("O204").Value = 1.117037037037037
$ws.Range("P204").Value = 4.7
$ws.Range("Q204").Value = 4.8
$ws.Range("R204").Value = 4.5
$ws.Range("S204").Value = 5
$ws.Range("T204").Value = 6.6
$ws.Range("U204").Value = 6.5
$ws.Range("V204").Value = 0.02973199381914737
$ws.Range("W204").Value = 0.03145225477868011
$ws.Range("X204").Value = 0.0311204687151446
$ws.Range("Y204").Value = 0.03032342634296706

# ---- Copy date-cell formatting (style index used by column A) down to the new rows ----
$ws.Range("A200").Copy()
$ws.Range("A201:A204").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

Write-Host "Edit complete"
